$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = "Record"
$ws.Range("B37").Value = "RJ Record"
$ws.Range("C37").Value = "Transporte"
$ws.Range("D37").Value = "2025-04-01T19:03"
$ws.Range("E37").Value = "Negativo"
$ws.Range("F37").Value = @'
Impasse em Campos. Motoristas de vans dizem que só voltam ao trabalho se prefeitura fizer repasse. Essa é a segunda paralisação do setor C, que faz a linha da Região Norte, que atende sete localidades. Moradores foram pegos de surpresa ontem no final da tarde. Entrevista com passageiros. Passageiros pagam R$ 3,50 e o valor é complementado pela prefeitura. Permissionários na frente da prefeitura. Eles querem falar com representante da prefeitura e prefeito Wladimir Garotinho. Eles alegam que teve o repasse ontem à tarde, mas não receberam o valor devido de acordo com a bilhetagem. Valor estaria com inconsistências. Receberam de 80% a 90% a menos do que deveriam receber. Pediram através do judiciário que atende aos permissionários para que fossem atendidos e saber o que realmente está acontecendo. Sistema de bilhetagem foi implantado a pedido da prefeitura em 2019.Alegam que não tem como trabalhar desse jeito. Entrevista com passageiros e com permissionário Jefferson Oliveira. *Com nota do IMTT* *matéria* Repórter *ao vivo*, às 19h06, atualiza as informações. Ainda não há consenso. Uma reunião estaria acontecendo à tarde. Tema também foi abordado no Balanço Geral.
'@ + [char]10

$ws.Rows(37).AutoFit()
